$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes existing data rows 2..26 down to 3..27).
$ws.Rows.Item(2).Insert(-4121, 0)

# The inserted row inherits formatting from the header row (bold) - strip it
# back to the plain/unstyled look used by the rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Column D uses a custom date-time display format throughout the data rows;
# restore it on the new row so the underlying serial number renders the same way.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the new data record.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44496
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112013
$ws.Range("G2").Value = "Alcachofa"
$ws.Range("H2").Value = "Madrigal"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("N2").Value = "$/caja 40 unidades"
$ws.Range("O2").Value = "Provincia del Elquí"
$ws.Range("P2").Value = 288
$ws.Range("Q2").Value = 40
$ws.Range("R2").Value = "Hortaliza"
